$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 updates
$ws.Range("D8").Value = 0.99
$ws.Range("G8").Value = 0.54
$ws.Range("J8").Value = 0.48
$ws.Range("M8").Value = 0.65
$ws.Range("P8").Value = 0.8100000000000001

# Row 9 updates
$ws.Range("D9").Value = 0.99
$ws.Range("G9").Value = 0.87
$ws.Range("J9").Value = 0.86
$ws.Range("M9").Value = 0.92
$ws.Range("P9").Value = 0.96
